# The edit removes the "pc" data row (row 23) entirely - deleting the row
# shifts every row below it up by one (so former row 24 "np_l" becomes the
# new row 23, etc.), and appends a brand-new row of data (labelled "zy_r")
# at the bottom of the table (new row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 23 ("pc"), shifting rows 24:48 up by one.
$ws.Rows("23:23").Delete()

# Append the new trailing row of data that the commit introduced.
$ws.Range("A48").Value = "zy_r"
$ws.Range("B48").Value = -58.78
$ws.Range("C48").Value = 3.67
$ws.Range("D48").Value = 66.66

# Reflect the final selection/scroll position left behind by the edit.
$excel.ActiveWindow.ScrollRow = 32
[void]$ws.Range("D48").Select()
